# Roboflow annotation report 7/24/2025
# Add a new weekly progress row (row 65) to the tracking table, mirroring
# the layout/format of the prior row (64), and resize the Excel table so
# the new row is included.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the structured table (Table1) to include the new row.
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("D4:J65"))

# Copy formatting (styles, number formats, borders, row height) from the
# last existing data row down into the new row.
$ws.Range("D64:J64").Copy()
$ws.Range("D65:J65").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D65:J65").RowHeight = $ws.Range("D64:J64").RowHeight

# Fill in the new row's data.
$ws.Range("D65").Value = "24/7/2028"
$ws.Range("E65").Value = 380
$ws.Range("F65").Value = 950
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 1012
$ws.Range("J65").Value = "N/A"

# Match the author's final cell selection.
$null = $ws.Range("F65").Select()
